$wb = $excel.ActiveWorkbook

# --- Sheet "include" (sheet1) ---
$wsInclude = $wb.Worksheets.Item("include")

$wsInclude.Range("A1").Value = "type"
$wsInclude.Range("B1").Value = "value"
$wsInclude.Range("A2").Value = "path"
$wsInclude.Range("B2").Value = "C:\Temp\images1\"
$wsInclude.Range("A3").Value = "path"
$wsInclude.Range("B3").Value = "C:\Temp\images2\images2-sub2\images2-sub2-sub2\"

$wsInclude.Range("E10").Select()

# --- Sheet "exclude" (sheet2) ---
$wsExclude = $wb.Worksheets.Item("exclude")

$wsExclude.Range("A1").Value = "type"
$wsExclude.Range("B1").Value = "value"
$wsExclude.Range("A2").Value = "path"
$wsExclude.Range("B2").Value = "C:\Temp\images\folder2\folder2-sub2"
$wsExclude.Range("A3").Value = "path"
$wsExclude.Range("B3").Value = "C:\Temp\images\folder2\folder2-sub1"
$wsExclude.Range("B4").Value = "C:\Temp\images2\images2-sub2\images2-sub2-sub2\images2-sub2-sub2-sub1"

$wsExclude.Range("B12").Select()
